$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "security/sha/runme_large.sh"
$ws.Range("B16").Value = 0.01
$ws.Range("C16").Value = 0.01
$ws.Range("D16").Value = 0

$ws.Range("A51").Select()
